# Regenerate save_data to use K instead of Strike#: update column G (K) values
# for rows 2-24 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 6
    4  = 2
    5  = 5
    6  = 4
    7  = 6
    8  = 6
    9  = 8
    10 = 6
    11 = 7
    12 = 7
    13 = 5
    14 = 3
    15 = 1
    16 = 5
    17 = 5
    18 = 7
    19 = 9
    20 = 5
    21 = 3
    22 = 6
    23 = 4
    24 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
